$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.096.66"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.01%  "
$ws.Range("D3").Value = "'1.867.75"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.21%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'307.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.98%  "
$ws.Range("E6").Value = "  +0.21%  "
$ws.Range("D7").Value = "'0.5094"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.49%  "
$ws.Range("D8").Value = "'0.3737"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.25%  "
$ws.Range("D9").Value = "'0.07153"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.35%  "
$ws.Range("D10").Value = "'0.8878"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.70%  "
$ws.Range("D11").Value = "'20.59"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.11%  "
$ws.Range("D12").Value = "'1.866.71"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.46%  "
$ws.Range("D13").Value = "'0.07545"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.67%  "
$ws.Range("D14").Value = "'5.318"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.12%  "
$ws.Range("D15").Value = "'89.31"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.91%  "
$ws.Range("D16").Value = "'1.000"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.15%  "
$ws.Range("D17").Value = "'0.000008465"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.31%  "
$ws.Range("D18").Value = "'14.12"
$ws.Range("D18").Style = "Normal"
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("D20").Value = "'27.142.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.92%  "
$ws.Range("D21").Value = "'5.065"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.35%  "
$ws.Range("D22").Value = "'2.105.16"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.71%  "
$ws.Range("D23").Value = "'10.57"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.72%  "
$ws.Range("E24").Value = "  -2.01%  "
$ws.Range("D25").Value = "'150.57"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.67%  "
$ws.Range("D26").Value = "'1.833"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.67%  "
$ws.Range("E27").Value = "  -2.78%  "
$ws.Range("D28").Value = "'2.098"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.05%  "
$ws.Range("D29").Value = "'112.53"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.60%  "
$ws.Range("D30").Value = "'4.748"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.78%  "
$ws.Range("D31").Value = "'4.685"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.77%  "
$ws.Range("D32").Value = "'0.09053"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.16%  "
$ws.Range("D33").Value = "'0.05129"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Value = "'3.097"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.52%  "
$ws.Range("D35").Value = "'1.161"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.30%  "
$ws.Range("D36").Value = "'0.7367"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.22%  "
$ws.Range("D37").Value = "'0.02039"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.29%  "
$ws.Range("D38").Value = "'2.485"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.23%  "
$ws.Range("D39").Value = "'3.042"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.78%  "
$ws.Range("E40").Value = "  -1.27%  "
$ws.Range("D41").Value = "'0.5339"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.99%  "
$ws.Range("D42").Value = "'6.606"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.12%  "
$ws.Range("D43").Value = "'115.91"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.02%  "
$ws.Range("D44").Value = "'8.347"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.08%  "
$ws.Range("D45").Value = "'0.1473"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.10%  "
$ws.Range("D46").Value = "'0.4640"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.16%  "
$ws.Range("D47").Value = "'1.001"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.24%  "
$ws.Range("D48").Value = "'10.03"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.28%  "
$ws.Range("D49").Value = "'1.565"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.64%  "
$ws.Range("D50").Value = "'64.56"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.64%  "
$ws.Range("D51").Value = "'36.42"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.16%  "
